# Electricity sector calibrations: increase RAF values to 0.9 for technologies
# less than that, on the "RAF-generation" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RAF-generation")

# hard coal: 0.6 -> 0.9
$ws.Range("B2").Value = 0.9

# biomass: 0.6 -> 0.9
$ws.Range("B10").Value = 0.9

# geothermal: 0.85 -> 0.9
$ws.Range("B11").Value = 0.9

# Update selection to match saved state
$ws.Range("B3").Select()
